{"js": "// Fix typo: \"Thank you for your consultation received on <Log Date>.....\"\n// should read \"...received on <Log Date>.\" \u2014 collapse the stray extra dots\n// (split across two separate runs) down to a single trailing period.\nconst body = context.document.body;\n\nconst results = body.search(\" <Log Date>.....\", { matchCase: true, matchWildcards: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '<Log Date>.....' text to fix.\");\n}\n\nresults.items[0].insertText(\" <Log Date>.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix typo: \"Thank you for your consultation received on <Log Date>.....\"\n# should read \"...received on <Log Date>.\" (collapse the stray extra dots\n# into a single period).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"<Log Date>.....\"\n$find.Replacement.Text = \"<Log Date>.\"\n$find.Forward = $true\n$find.Wrap = 0            # wdFindStop - don't wrap past the single intended match\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n    [ref]$find.Text,\n    [ref]$find.MatchCase,\n    [ref]$find.MatchWholeWord,\n    $null,\n    $null,\n    $null,\n    [ref]$find.Forward,\n    [ref]$find.Wrap,\n    [ref]$find.Format,\n    [ref]$find.Replacement.Text,\n    2                       # wdReplaceAll\n) | Out-Null\n"}
